$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure new rows (135-138) inherit the same cell style/format as the
# existing data rows (s="1": thin border, centered) before writing values.
$ws.Range("A134").Copy()
$ws.Range("A135:A138").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$data = @(
    @(2, "Sample_0100", 0),
    @(3, "Sample_0226", 1),
    @(4, "Sample_0085", 1),
    @(5, "Sample_0047", 1),
    @(6, "Sample_0154", 0),
    @(7, "Sample_0218", 1),
    @(8, "Sample_0088", 1),
    @(9, "Sample_0251", 1),
    @(10, "Sample_0160", 1),
    @(11, "Sample_0038", 0),
    @(12, "Sample_0206", 0),
    @(13, "Sample_0208", 0),
    @(14, "Sample_0242", 1),
    @(15, "Sample_0006", 0),
    @(16, "Sample_0093", 1),
    @(17, "Sample_0163", 1),
    @(18, "Sample_0204", 0),
    @(19, "Sample_0219", 1),
    @(20, "Sample_0074", 0),
    @(21, "Sample_0246", 1),
    @(22, "Sample_0033", 1),
    @(23, "Sample_0187", 1),
    @(24, "Sample_0148", 1),
    @(25, "Sample_0053", 1),
    @(26, "Sample_0018", 0),
    @(27, "Sample_0185", 1),
    @(28, "Sample_0215", 1),
    @(29, "Sample_0013", 1),
    @(30, "Sample_0039", 1),
    @(31, "Sample_0168", 0),
    @(32, "Sample_0011", 1),
    @(33, "Sample_0052", 1),
    @(34, "Sample_0253", 1),
    @(35, "Sample_0201", 0),
    @(36, "Sample_0157", 1),
    @(37, "Sample_0200", 0),
    @(38, "Sample_0017", 0),
    @(39, "Sample_0115", 0),
    @(40, "Sample_0159", 1),
    @(41, "Sample_0202", 0),
    @(42, "Sample_0003", 1),
    @(43, "Sample_0235", 1),
    @(44, "Sample_0238", 1),
    @(45, "Sample_0124", 0),
    @(46, "Sample_0170", 0),
    @(47, "Sample_0117", 1),
    @(48, "Sample_0091", 1),
    @(49, "Sample_0186", 1),
    @(50, "Sample_0040", 0),
    @(51, "Sample_0087", 1),
    @(52, "Sample_0176", 0),
    @(53, "Sample_0240", 1),
    @(54, "Sample_0212", 1),
    @(55, "Sample_0037", 1),
    @(56, "Sample_0034", 1),
    @(57, "Sample_0207", 0),
    @(58, "Sample_0203", 0),
    @(59, "Sample_0139", 1),
    @(60, "Sample_0167", 0),
    @(61, "Sample_0137", 0),
    @(62, "Sample_0142", 0),
    @(63, "Sample_0086", 1),
    @(64, "Sample_0112", 1),
    @(65, "Sample_0180", 1),
    @(66, "Sample_0144", 1),
    @(67, "Sample_0135", 0),
    @(68, "Sample_0250", 1),
    @(69, "Sample_0224", 1),
    @(70, "Sample_0094", 1),
    @(71, "Sample_0152", 0),
    @(72, "Sample_0210", 0),
    @(73, "Sample_0090", 0),
    @(74, "Sample_0234", 1),
    @(75, "Sample_0009", 1),
    @(76, "Sample_0164", 1),
    @(77, "Sample_0103", 0),
    @(78, "Sample_0193", 1),
    @(79, "Sample_0236", 1),
    @(80, "Sample_0213", 1),
    @(81, "Sample_0241", 1),
    @(82, "Sample_0097", 1),
    @(83, "Sample_0138", 1),
    @(84, "Sample_0057", 1),
    @(85, "Sample_0014", 0),
    @(86, "Sample_0169", 0),
    @(87, "Sample_0149", 0),
    @(88, "Sample_0141", 0),
    @(89, "Sample_0232", 1),
    @(90, "Sample_0237", 1),
    @(91, "Sample_0147", 1),
    @(92, "Sample_0015", 1),
    @(93, "Sample_0020", 1),
    @(94, "Sample_0110", 0),
    @(95, "Sample_0227", 1),
    @(96, "Sample_0191", 1),
    @(97, "Sample_0025", 1),
    @(98, "Sample_0126", 0),
    @(99, "Sample_0231", 1),
    @(100, "Sample_0194", 1),
    @(101, "Sample_0223", 0),
    @(102, "Sample_0254", 1),
    @(103, "Sample_0116", 0),
    @(104, "Sample_0105", 1),
    @(105, "Sample_0166", 0),
    @(106, "Sample_0099", 0),
    @(107, "Sample_0252", 1),
    @(108, "Sample_0109", 1),
    @(109, "Sample_0095", 0),
    @(110, "Sample_0146", 1),
    @(111, "Sample_0032", 1),
    @(112, "Sample_0108", 1),
    @(113, "Sample_0220", 1),
    @(114, "Sample_0249", 1),
    @(115, "Sample_0189", 1),
    @(116, "Sample_0083", 0),
    @(117, "Sample_0217", 1),
    @(118, "Sample_0244", 1),
    @(119, "Sample_0004", 1),
    @(120, "Sample_0101", 0),
    @(121, "Sample_0209", 0),
    @(122, "Sample_0190", 1),
    @(123, "Sample_0211", 0),
    @(124, "Sample_0162", 1),
    @(125, "Sample_0158", 1),
    @(126, "Sample_0175", 0),
    @(127, "Sample_0145", 0),
    @(128, "Sample_0131", 1),
    @(129, "Sample_0181", 1),
    @(130, "Sample_0255", 1),
    @(131, "Sample_0031", 1),
    @(132, "Sample_0248", 1),
    @(133, "Sample_0125", 0),
    @(134, "Sample_0165", 0),
    @(135, "Sample_0022", 1),
    @(136, "Sample_0150", 0),
    @(137, "Sample_0239", 1),
    @(138, "Sample_0161", 1)
)

foreach ($row in $data) {
    $r = $row[0]
    $sampleId = $row[1]
    $label = $row[2]
    $ws.Cells.Item($r, 1).Value = $sampleId
    $ws.Cells.Item($r, 2).Value = $label
}
